$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("seats")

# New data rows (user_id, loadable_place_id, timestamp) for rows 2-6
$data = @(
    @(0, "bogdan@nu.edu.kz",          "Test event 1; Cinema hall, floor 1, upper middle section, row F, seat 5", 45509.62525357998),
    @(1, "admin",                     "Test event 1; Cinema hall, floor 1, upper middle section, row H, seat 3", 45508.49005549768),
    @(2, "admin",                     "Test event 1; Cinema hall, floor 1, upper middle section, row D, seat 2", 45508.49010777778),
    @(3, "bogdan.yakupov@nu.edu.kz",  "Test event 1; Cinema hall, floor 1, upper left section, row F, seat 3",  45519.62325730543),
    @(4, "adema.akizhanova@nu.edu.kz","Test event 1; Cinema hall, floor 1, upper left section, row G, seat 6",  45519.62876581898)
)

# Clear rows 7-9 (previously used, now removed)
$ws.Range("A7:D9").Clear()

# Write column by column so new shared-string entries are interned
# in the same order the source system produced them (all user_id
# values, then all loadable_place_id values).
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
